$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.375.94'
$ws.Range('E2').Value = '  +3.81%  '
$ws.Range('D3').Value = '3.074.43'
$ws.Range('E3').Value = '  +5.79%  '
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').Value = '514.44'
$ws.Range('E5').Value = '  +3.39%  '
$ws.Range('D6').Value = '141.70'
$ws.Range('E6').Value = '  +7.13%  '
$ws.Range('E7').Value = '  -0.10%  '
$ws.Range('D8').Value = '0.435'
$ws.Range('E8').Value = '  +4.05%  '
$ws.Range('E9').Value = '  +2.04%  '
$ws.Range('D10').Value = '0.109'
$ws.Range('E10').Value = '  +5.12%  '
$ws.Range('D11').Value = '0.373'
$ws.Range('E11').Value = '  +7.64%  '
$ws.Range('D12').Value = '3.596.14'
$ws.Range('E12').Value = '  +5.50%  '
$ws.Range('E13').Value = '  +3.07%  '
$ws.Range('D14').Value = '25.57'
$ws.Range('E14').Value = '  +0.01%  '
$ws.Range('E15').Value = '  +4.15%  '
$ws.Range('D16').Value = '57.426.32'
$ws.Range('E16').Value = '  +3.78%  '
$ws.Range('D17').Value = '3.071.28'
$ws.Range('E17').Value = '  +5.55%  '
$ws.Range('D18').Value = '6.05'
$ws.Range('E18').Value = '  +2.07%  '
$ws.Range('D19').Value = '13.03'
$ws.Range('E19').Value = '  +4.45%  '
$ws.Range('D20').Value = '8.17'
$ws.Range('E20').Value = '  +6.97%  '
$ws.Range('D21').Value = '339.20'
$ws.Range('E21').Value = '  +8.75%  '
$ws.Range('E22').Value = '  +0.21%  '
$ws.Range('D23').Value = '0.500'
$ws.Range('E23').Value = '  +3.60%  '
$ws.Range('D24').Value = '65.62'
$ws.Range('E24').Value = '  +5.35%  '
$ws.Range('D25').Value = '0.173'
$ws.Range('E25').Value = '  +8.58%  '
$ws.Range('E26').Value = '  +0.32%  '
$ws.Range('D27').Value = '0.0₃0939'
$ws.Range('E27').Value = '  +13.68%  '
$ws.Range('D28').Value = '6.46'
$ws.Range('E28').Value = '  +2.67%  '
$ws.Range('D29').Value = '7.15'
$ws.Range('E29').Value = '  +4.50%  '
$ws.Range('D30').Value = '1.81'
$ws.Range('E30').Value = '  +4.02%  '
$ws.Range('D31').Value = '20.81'
$ws.Range('E31').Value = '  +5.85%  '
$ws.Range('E32').Value = '  +5.34%  '
$ws.Range('D33').Value = '154.00'
$ws.Range('E33').Value = '  +2.60%  '
$ws.Range('D34').Value = '4.55'
$ws.Range('E34').Value = '  +5.51%  '
$ws.Range('D35').Value = '5.90'
$ws.Range('E35').Value = '  +6.65%  '
$ws.Range('D36').Value = '26.21'
$ws.Range('E36').Value = '  +9.50%  '
$ws.Range('E37').Value = '  +5.53%  '
$ws.Range('D38').Value = '0.0679'
$ws.Range('E38').Value = '  +6.22%  '
$ws.Range('D39').Value = '3.110.77'
$ws.Range('E39').Value = '  +5.89%  '
$ws.Range('D40').Value = '37.03'
$ws.Range('E40').Value = '  +2.32%  '
$ws.Range('D41').Value = '3.86'
$ws.Range('E41').Value = '  +5.59%  '
$ws.Range('D42').Value = '0.671'
$ws.Range('E42').Value = '  +6.03%  '
$ws.Range('E43').Value = '  -0.22%  '
$ws.Range('D44').Value = '2.263.09'
$ws.Range('E44').Value = '  +7.71%  '
$ws.Range('B45').Value = 'Stacks'
$ws.Range('C45').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D45').Value = '1.39'
$ws.Range('E45').Value = '  +5.38%  '
$ws.Range('B46').Value = 'VeChain'
$ws.Range('C46').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D46').Value = '0.0251'
$ws.Range('E46').Value = '  +9.05%  '
$ws.Range('D47').Value = '0.956'
$ws.Range('E47').Value = '  +5.42%  '
$ws.Range('D48').Value = '20.20'
$ws.Range('E48').Value = '  +9.60%  '
$ws.Range('D49').Value = '5.87'
$ws.Range('E49').Value = '  -0.79%  '
$ws.Range('D50').Value = '0.0869'
$ws.Range('E50').Value = '  +4.40%  '
$ws.Range('B51').Value = 'TheGraph'
$ws.Range('C51').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D51').Value = '0.181'
$ws.Range('E51').Value = '  +6.03%  '
